$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.907.89'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.635.28'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.91'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5051'
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2568'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06351'
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.65'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07736'
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.271'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.637.61'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5427'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅7724'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.03'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.931.57'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.425'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.82'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.901'
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.096'
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.887'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.80'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1237'
$ws.Range('E26').Value = '  +7.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.802'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.55'
$ws.Range('E28').Value = '  -1.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.235'
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04864'
$ws.Range('E30').Value = '  -3.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.240'
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.191'
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.546'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9093'
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.572'
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.122.72'
$ws.Range('E37').Value = '  -1.24%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5488'
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01556'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.573'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8018'
$ws.Range('E42').Value = '  -2.30%  '
$ws.Range('E43').Value = '  -8.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.43'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.770.45'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4475'
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.92'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05158'
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.509'
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9998'
$ws.Range('E51').Value = '  -0.69%  '